# Update the four "vegas" game rows (away team, spread, total) with new
# matchups/lines. All other cells on the sheet (E:J formulas, B9:C12 summary
# strings, the G2:G5 / J2:J5 what-if Data Tables, and the L/S helper columns)
# recompute automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LAC @ -6/40.5  ->  PIT @ -6/30
$ws.Range("B2").Value = "PIT"
$ws.Range("D2").Value = 30

# Row 3: KC @ -6.5/42  ->  KC @ -2.5/47
$ws.Range("C3").Value = -2.5
$ws.Range("D3").Value = 47

# Row 4: was blank  ->  MIN @ -3/40.1
$ws.Range("B4").Value = "MIN"
$ws.Range("C4").Value = -3
$ws.Range("D4").Value = 40.1

# Row 5: LAR @ -4/39.5  ->  BAL @ -7/42.5
$ws.Range("B5").Value = "BAL"
$ws.Range("C5").Value = -7
$ws.Range("D5").Value = 42.5

# Recalculate so the Data Tables (G2:G5, J2:J5) and all dependent formulas
# pick up the new inputs.
$excel.Calculate()

# Selection moved to D5 (matches the last-edited input cell).
[void]$ws.Range("D5").Select()
